$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The student IDs in column A (rows 1-70) were re-mapped from 7001..7070
# down to plain sequential integers 1..70 (student id -> int).
for ($i = 1; $i -le 70; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}

# Reset the view/selection away from the old "A57 / E67" spot to A1:A70,
# matching the saved sheet view in the updated workbook.
$ws.Range("A1:A70").Select()

$wb.Save()
